$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Huwebes"
$ws.Range("B1").Value = "DATE"

$ws.Range("A2").Value = "Cavite"
$ws.Range("B2").Value = "LOCATION"

$ws.Range("A3").Value = "DOJ"
$ws.Range("B3").Value = "LOCATION"

$ws.Range("A4").Value = "NBI"
$ws.Range("B4").Value = "LOCATION"
$ws.Range("D4").Value = "Person"
$ws.Range("E4").Value = "Location"
$ws.Range("F4").Value = "Date"

$ws.Range("A5").Value = "Dacer at Corbito"
$ws.Range("B5").Value = "PERSON"
$ws.Range("D5").Value = "Alex Diloy"
$ws.Range("E5").Value = "Cavite"
$ws.Range("F5").Value = "Huwebes"

$ws.Range("A6").Value = "Dacer"
$ws.Range("B6").Value = "PERSON"
$ws.Range("D6").Value = "Corbito"

$ws.Range("A7").Value = "Ebdane"
$ws.Range("B7").Value = "PERSON"
$ws.Range("D7").Value = "Dacer"

$ws.Range("A8").Value = "Emmanuel Corbito"
$ws.Range("B8").Value = "PERSON"
$ws.Range("D8").Value = "Diloy"

$ws.Range("A9").Value = "Grace Amargo at Joy Cantos"
$ws.Range("B9").Value = "PERSON"
$ws.Range("D9").Value = "Ebdane"

$ws.Range("A10").Value = "Jimmy Lopez at Alex Diloy,"
$ws.Range("B10").Value = "PERSON"
$ws.Range("D10").Value = "Emmanuel Corbito"

$ws.Range("A11").Value = "Lopez at Diloy"
$ws.Range("B11").Value = "PERSON"
$ws.Range("D11").Value = "Grace Amargo"

$ws.Range("A12").Value = "NBI Director Reynaldo Wycoco"
$ws.Range("B12").Value = "PERSON"
$ws.Range("D12").Value = "Jimmy Lopez"

$ws.Range("A13").Value = "NBI-National Capital Region Director Samuel Ong"
$ws.Range("B13").Value = "PERSON"
$ws.Range("D13").Value = "Joy Cantos"

$ws.Range("A14").Value = "Ong"
$ws.Range("B14").Value = "PERSON"
$ws.Range("D14").Value = "Lopez"

$ws.Range("A15").Value = "PAOCTF Chief Director Hermogenes Ebdane"
$ws.Range("B15").Value = "PERSON"
$ws.Range("D15").Value = "NBI Director Reynaldo Wycoco"

$ws.Range("A16").Value = "PAOCTF-Visayas Chief Sr. Supt. Teofilo Vina"
$ws.Range("B16").Value = "PERSON"
$ws.Range("D16").Value = "NBI-National Capital Region Director Samuel Ong"

$ws.Range("A17").Value = "Pangulong Gloria Macapagal-Arroyo"
$ws.Range("B17").Value = "PERSON"
$ws.Range("D17").Value = "Ong"

$ws.Range("A18").Value = "Salvador Bubby`" Dacer`""
$ws.Range("B18").Value = "PERSON"
$ws.Range("D18").Value = "Pangulong Gloria Macapagal-Arroyo"

$ws.Range("D19").Value = "PAOCTF Chief Director Hermogenes Ebdane"

$ws.Range("D20").Value = "PAOCTF-Visayas Chief Sr. Supt. Teofilo Vina"

$ws.Range("D21").Value = "Salvador `"Bubby`" Dacer"

$ws.Range("A19").ClearContents()
$ws.Range("B19").ClearContents()
$ws.Range("A20").ClearContents()
$ws.Range("B20").ClearContents()
$ws.Range("A21").ClearContents()
$ws.Range("B21").ClearContents()

$ws.Range("D5").Select()
